$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings: force text format
# so Excel keeps them as text (preserving exact formatting/trailing zeros)
# rather than silently converting them to floating point numbers.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = "67.017.58"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "2.524.47"
$ws.Range("E3").Value = "  -2.32%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "589.51"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("D6").Value = "172.87"
$ws.Range("E6").Value = "  +3.98%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.527"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "2.524.64"
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("D10").Value = "0.138"
$ws.Range("E10").Value = "  +0.21%  "
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").Value = "5.15"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("E13").Value = "  -3.89%  "
$ws.Range("D14").Value = "26.53"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "2.985.97"
$ws.Range("E15").Value = "  -2.29%  "
$ws.Range("D16").Value = "0.0000176"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").Value = "66.832.93"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "2.523.98"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "8.08"
$ws.Range("E19").Value = "  +4.71%  "
$ws.Range("D20").Value = "11.32"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "353.79"
$ws.Range("E21").Value = "  +0.46%  "
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("D23").Value = "4.62"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "1.99"
$ws.Range("E24").Value = "  +5.64%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "69.62"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").Value = "9.93"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").Value = "0.0₃0975"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").Value = "531.84"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "8.14"
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("D33").Value = "1.33"
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("E35").Value = "  -0.98%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").Value = "157.80"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "18.60"
$ws.Range("E39").Value = "  -0.77%  "
$ws.Range("D40").Value = "18.44"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "5.13"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").Value = "2.49"
$ws.Range("E45").Value = "  +3.09%  "
$ws.Range("D46").Value = "148.98"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D48").Value = "0.0₆0277"
$ws.Range("E48").Value = "  -3.01%  "
$ws.Range("D49").Value = "3.69"
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("E51").Value = "  -0.12%  "
